# working on logging bug
# Append new diary rows (51-60) to the "Arbeitszeit" log sheet, mirroring
# the existing E:I layout (Datum / Zeit / Einheit / Tätigkeit / ...).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 51
$ws.Range("E51").Value2 = 43733
$ws.Range("E51").NumberFormat = "m/d/yy"
$ws.Range("F51").Value = 4
$ws.Range("G51").Value = "Stunden"
$ws.Range("H51").Value = "Dokumentation"
$ws.Range("I51").Value = "Pflichtenheft, Lift-Pitch"

# Row 52
$ws.Range("E52").Value2 = 43736
$ws.Range("E52").NumberFormat = "m/d/yy"
$ws.Range("F52").Value = 4
$ws.Range("G52").Value = "Stunden"
$ws.Range("H52").Value = "Pflichtenheft"

# Row 53
$ws.Range("E53").Value2 = 43739
$ws.Range("E53").NumberFormat = "m/d/yy"
$ws.Range("F53").Value = 3
$ws.Range("G53").Value = "Stunden"
$ws.Range("H53").Value = "Pflichtenheft, Exposé"

# Row 54
$ws.Range("E54").Value2 = 43742
$ws.Range("E54").NumberFormat = "m/d/yy"
$ws.Range("F54").Value = 2
$ws.Range("G54").Value = "Stunden"
$ws.Range("H54").Value = "Fertigstellung Pflichtenheft"

# Row 55
$ws.Range("E55").Value2 = 43759
$ws.Range("E55").NumberFormat = "m/d/yy"
$ws.Range("F55").Value = 2
$ws.Range("G55").Value = "Stunden"
$ws.Range("H55").Value = "Ausarbeitung Dokumentation"

# Row 56
$ws.Range("E56").Value2 = 43766
$ws.Range("E56").NumberFormat = "m/d/yy"
$ws.Range("F56").Value = 3
$ws.Range("G56").Value = "Stunden"
$ws.Range("H56").Value = "Inhaltsangabe Präsentation"

# Row 57
$ws.Range("E57").Value2 = 43773
$ws.Range("E57").NumberFormat = "m/d/yy"
$ws.Range("F57").Value = 2
$ws.Range("G57").Value = "Stunden"
$ws.Range("H57").Value = "Ausarbeitung Dokumentation"

# Row 58
$ws.Range("E58").Value2 = 43777
$ws.Range("E58").NumberFormat = "m/d/yy"
$ws.Range("F58").Value = 1
$ws.Range("G58").Value = "Stunden"
$ws.Range("H58").Value = "Dokumenation"

# Row 59
$ws.Range("E59").Value2 = 43780
$ws.Range("E59").NumberFormat = "m/d/yy"
$ws.Range("F59").Value = 2
$ws.Range("G59").Value = "Stunden"
$ws.Range("H59").Value = "Layout Präsentation"

# Row 60
$ws.Range("E60").Value2 = 43782
$ws.Range("E60").NumberFormat = "m/d/yy"
$ws.Range("F60").Value = 1
$ws.Range("G60").Value = "Stunden"
$ws.Range("H60").Value = "Dokumentation Beifügungen"

# Move the selection to mirror the author's final cursor position after
# adding the new entries (view scrolls down as rows are appended).
$ws.Range("F61").Select()
